$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains four OpenL "Method" tables (each a 2-row block: a header
# row followed by a body row) plus an "Environment" table at the bottom.
# This change adds a "properties" sub-header row (properties | version |
# 0.0.2) directly under each Method table's header row, pushing every row
# below it down by one. We insert the four new rows from the bottom of the
# sheet upward so that the original (pre-edit) row numbers for the
# not-yet-processed inserts stay valid.

$ws.Rows(16).Insert()
$ws.Rows(12).Insert()
$ws.Rows(8).Insert()
$ws.Rows(4).Insert()

# Give the new rows the same bordered / centered look as the table header
# cells right above them (copy format only, keep the new cells unmerged).
$ws.Range("B3").Copy()
$ws.Range("B4:D4").PasteSpecial(-4122)
$ws.Range("B9:D9").PasteSpecial(-4122)
$ws.Range("B14:D14").PasteSpecial(-4122)
$ws.Range("B19:D19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new "properties" rows.
$ws.Range("B4").Value2 = "properties"
$ws.Range("C4").Value2 = "version"
$ws.Range("D4").Value2 = "0.0.2"

$ws.Range("B9").Value2 = "properties"
$ws.Range("C9").Value2 = "version"
$ws.Range("D9").Value2 = "0.0.2"

$ws.Range("B14").Value2 = "properties"
$ws.Range("C14").Value2 = "version"
$ws.Range("D14").Value2 = "0.0.2"

$ws.Range("B19").Value2 = "properties"
$ws.Range("C19").Value2 = "version"
$ws.Range("D19").Value2 = "0.0.2"

# Restore the active cell/selection used in the edited workbook.
$ws.Range("F16").Select() | Out-Null
